$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 446; this shifts the existing rows
# 446:554 down to 447:555 (carrying their values/formatting along),
# leaving row 446 blank (inheriting the D-column date style from the
# row above because Excel's default insert picks up formatting from
# above the insertion point).
$ws.Rows(446).Insert()

# Populate the newly inserted row 446 with the new weekly data point.
# Columns A,B,C,E,F,G,H,I,J,K,Q,T are constant across every data row
# in this sheet, so we just repeat the same values used everywhere
# else. L (Calidad) and R (Origen) happen to match what the row
# directly below already has (Primera / Perú), M/N/O/P/S are the new
# reported volume & price figures, and D is the new report date.
$ws.Range("A446").Value = 3
$ws.Range("B446").Value = "Femacal de La Calera"
$ws.Range("C446").Value = "Coquimbo"
$ws.Range("D446").Value2 = 44943
$ws.Range("E446").Value = 5
$ws.Range("F446").Value = "Fruta"
$ws.Range("G446").Value = 100108
$ws.Range("H446").Value = "Tropicales y subtropicales"
$ws.Range("I446").Value = 100108002
$ws.Range("J446").Value = "Mango"
$ws.Range("K446").Value = "Sin especificar"
$ws.Range("L446").Value = "Primera"
$ws.Range("M446").Value = 228
$ws.Range("N446").Value = 7000
$ws.Range("O446").Value = 7000
$ws.Range("P446").Value = 7000
$ws.Range("Q446").Value = '$/bandeja 4 kilos'
$ws.Range("R446").Value = "Perú"
$ws.Range("S446").Value = 1750
$ws.Range("T446").Value = 4
